# Remove the "Ver no Jupiter..." / copyright footer block that trails the
# "Requisitos" section (these paragraphs are no longer present in the
# rebuilt site output), while leaving the rest of the document untouched.
#
# Target region (inclusive), originally right after the requirement line
# "LOB1038: Física Experimental I (Requisito fraco)":
#   - (empty paragraph)
#   - "Ver no Jupiter Salvar em pdf Salvar em docx"
#   - "© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github
#      pages. Original theme under Creative Commons Attribution"

$d = $word.ActiveDocument

$count = $d.Paragraphs.Count

$anchorIndex = -1
for ($i = 1; $i -le $count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "LOB1038:*Requisito fraco*") {
        $anchorIndex = $i
        break
    }
}

$endIndex = -1
for ($i = 1; $i -le $count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "*Creative Commons Attribution*") {
        $endIndex = $i
        break
    }
}

if ($anchorIndex -gt 0 -and $endIndex -ge $anchorIndex) {
    $startPara = $d.Paragraphs($anchorIndex + 1)
    $endPara = $d.Paragraphs($endIndex)

    $rng = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $rng.Delete()
}
